# Weekly Fruta/hortaliza update: insert a new price-report row at row 75
# (Femacal de La Calera - Arándano (blue), Provincia de Curicó, fecha 2021-12-09),
# pushing the existing rows 75..166 down to 76..167.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 75; everything below shifts down one row
# and the sheet's used-range / dimension grows from T166 to T167 automatically.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly record.
$ws.Cells.Item(75, 1).Value = 3
$ws.Cells.Item(75, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).Value = 44539
$ws.Cells.Item(75, 5).Value = 5
$ws.Cells.Item(75, 6).Value = "Fruta"
$ws.Cells.Item(75, 7).Value = 100101
$ws.Cells.Item(75, 8).Value = "Berries"
$ws.Cells.Item(75, 9).Value = 100101001
$ws.Cells.Item(75, 10).Value = "Arándano (blue)"
$ws.Cells.Item(75, 11).Value = "Sin especificar"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 50
$ws.Cells.Item(75, 14).Value = 6000
$ws.Cells.Item(75, 15).Value = 6000
$ws.Cells.Item(75, 16).Value = 6000
$ws.Cells.Item(75, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(75, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(75, 19).Value = 3000
$ws.Cells.Item(75, 20).Value = 2
